$d = $word.ActiveDocument

function Split-NumberRun($searchText, $suffix) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    # Clear the matched run's text, then inject two fresh runs in its place
    # (kept inside the same <w:pPr> so the paragraph's style survives).
    $rng.Text = ""
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:pPr><w:pStyle w:val="style0"/></w:pPr>' +
        '<w:r><w:rPr/><w:t xml:space="preserve">' + $searchText + ' </w:t></w:r>' +
        '<w:r><w:rPr/><w:t>' + $suffix + '</w:t></w:r>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$rng.InsertXML($xml)
}

Split-NumberRun "36.86" "- with"
Split-NumberRun "45.47" "- without"

Split-NumberRun "19.02" "- with"
Split-NumberRun "19.67" "- without"

Split-NumberRun "1.04.74" "- with"
Split-NumberRun "1.07.03" "- without"

Split-NumberRun "25.96" "- with"
Split-NumberRun "16.05" "- without"

Split-NumberRun "15.54" "- with"
Split-NumberRun "15.93" "- without"

Split-NumberRun "40.10" "- with"
Split-NumberRun "47.97" "- without"
